$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - Bibi Cell Mundi
$ws.Range("H2").Value = 16440.91
$ws.Range("I2").Value = 14810.14
$ws.Range("J2").Value = 16474.68
$ws.Range("K2").Value = 1687
$ws.Range("AG2").Value = 102210.91

# Row 3 - Bibi Cell Vieiralves
$ws.Range("G3").Value = 12661.9
$ws.Range("H3").Value = 7883
$ws.Range("I3").Value = 7395
$ws.Range("J3").Value = 4342
$ws.Range("K3").Value = 6730
$ws.Range("AG3").Value = 73319.8

# Row 4 - Bibi Cell Manauara
$ws.Range("G4").Value = 2487
$ws.Range("H4").Value = 4148
$ws.Range("I4").Value = 3419.9
$ws.Range("J4").Value = 2395.49
$ws.Range("K4").Value = 4368
$ws.Range("L4").Value = 2742
$ws.Range("AG4").Value = 37276.89

# Row 5 - Bibi Cell Ponta Negra
$ws.Range("G5").Value = 2150.81
$ws.Range("H5").Value = 2966
$ws.Range("I5").Value = 2107.11
$ws.Range("J5").Value = 2296
$ws.Range("K5").Value = 2017.01
$ws.Range("L5").Value = 6974.9
$ws.Range("AG5").Value = 33509.11

# Row 6 - total
$ws.Range("G6").Value = 28248.43
$ws.Range("H6").Value = 31437.91
$ws.Range("I6").Value = 27732.15
$ws.Range("J6").Value = 25508.17
$ws.Range("K6").Value = 14802.01
$ws.Range("L6").Value = 9716.9
$ws.Range("AG6").Value = 246316.71
